# Appended data to 'Sheet1'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11. The migration-date column holds plain text (e.g. "2025-10-17")
# in every existing row, not a real date, so force text with a leading
# apostrophe (Excel's quote-prefix) and then reset the style back to
# "Normal" so no stray number-format style sticks to the cell.
$ws.Range("A11").Value = "'2025-10-23"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = "QQQ"
$ws.Range("C11").Value = "333OOO444"
$ws.Range("D11").Value = "Perungudi"

# E11 stays empty (same as E3:E9) but the cell itself still needs to exist
# as an empty text cell, so write the empty-text quote-prefix marker and
# reset the style the same way.
$ws.Range("E11").Value = "'"
$ws.Range("E11").Style = "Normal"
